# Remove "bhp" (column N) and "mileage" (column P) columns from the
# bulk-upload sample template. Deleting the entire columns shifts the
# trailing columns (category, images, storeId, discountPercent) left,
# matching the target layout, and also updates the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "mileage" column first (it is further right), then "bhp",
# so column letters for the earlier deletion stay valid.
$ws.Range("P1:P2").EntireColumn.Delete()
$ws.Range("N1:N2").EntireColumn.Delete()

# Reflect the new selection recorded in the saved file (O1, which is
# now "images" after the shift).
$ws.Range("O1").Select()
